{"js": "// Update the date line and all the two-digit multiplication problems.\nconst replacements = [\n  [\"2024-04-06 Saturday\", \"2024-04-07 Sunday\"],\n  [\"63\u00d742=\", \"25\u00d746=\"],\n  [\"38\u00d747=\", \"69\u00d748=\"],\n  [\"67\u00d757=\", \"72\u00d723=\"],\n  [\"98\u00d727=\", \"29\u00d720=\"],\n  [\"97\u00d798=\", \"18\u00d747=\"],\n  [\"24\u00d754=\", \"19\u00d749=\"],\n  [\"40\u00d740=\", \"99\u00d726=\"],\n  [\"23\u00d753=\", \"42\u00d787=\"],\n  [\"87\u00d772=\", \"56\u00d774=\"],\n  [\"23\u00d786=\", \"26\u00d795=\"],\n  [\"62\u00d734=\", \"64\u00d798=\"],\n  [\"51\u00d715=\", \"62\u00d732=\"],\n  [\"54\u00d748=\", \"37\u00d754=\"],\n  [\"29\u00d784=\", \"76\u00d775=\"],\n  [\"63\u00d755=\", \"90\u00d727=\"],\n  [\"36\u00d789=\", \"37\u00d742=\"],\n  [\"76\u00d773=\", \"30\u00d738=\"],\n  [\"70\u00d747=\", \"61\u00d779=\"],\n  [\"89\u00d732=\", \"14\u00d752=\"],\n  [\"73\u00d771=\", \"57\u00d715=\"],\n  [\"47\u00d783=\", \"83\u00d735=\"],\n  [\"38\u00d775=\", \"88\u00d772=\"],\n  [\"69\u00d714=\", \"18\u00d778=\"],\n  [\"35\u00d759=\", \"42\u00d718=\"],\n  [\"76\u00d758=\", \"13\u00d736=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"2024-04-06 Saturday\"; new = \"2024-04-07 Sunday\"},\n    @{old = \"63\u00d742=\"; new = \"25\u00d746=\"},\n    @{old = \"38\u00d747=\"; new = \"69\u00d748=\"},\n    @{old = \"67\u00d757=\"; new = \"72\u00d723=\"},\n    @{old = \"98\u00d727=\"; new = \"29\u00d720=\"},\n    @{old = \"97\u00d798=\"; new = \"18\u00d747=\"},\n    @{old = \"24\u00d754=\"; new = \"19\u00d749=\"},\n    @{old = \"40\u00d740=\"; new = \"99\u00d726=\"},\n    @{old = \"23\u00d753=\"; new = \"42\u00d787=\"},\n    @{old = \"87\u00d772=\"; new = \"56\u00d774=\"},\n    @{old = \"23\u00d786=\"; new = \"26\u00d795=\"},\n    @{old = \"62\u00d734=\"; new = \"64\u00d798=\"},\n    @{old = \"51\u00d715=\"; new = \"62\u00d732=\"},\n    @{old = \"54\u00d748=\"; new = \"37\u00d754=\"},\n    @{old = \"29\u00d784=\"; new = \"76\u00d775=\"},\n    @{old = \"63\u00d755=\"; new = \"90\u00d727=\"},\n    @{old = \"36\u00d789=\"; new = \"37\u00d742=\"},\n    @{old = \"76\u00d773=\"; new = \"30\u00d738=\"},\n    @{old = \"70\u00d747=\"; new = \"61\u00d779=\"},\n    @{old = \"89\u00d732=\"; new = \"14\u00d752=\"},\n    @{old = \"73\u00d771=\"; new = \"57\u00d715=\"},\n    @{old = \"47\u00d783=\"; new = \"83\u00d735=\"},\n    @{old = \"38\u00d775=\"; new = \"88\u00d772=\"},\n    @{old = \"69\u00d714=\"; new = \"18\u00d778=\"},\n    @{old = \"35\u00d759=\"; new = \"42\u00d718=\"},\n    @{old = \"76\u00d758=\"; new = \"13\u00d736=\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $r.new, 2)\n}\n"}
